$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("A1").Value = "garbageinfo@gmail.com"
$ws.Range("A2").Value = "ayahya.testitg@gmail.com"

# New rows
$ws.Range("A3").Value = "hello"
$ws.Range("A4").Value = "046pre056ty"

# Hyperlinks: replace A2's old hyperlink, add A1's and A2's new ones
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A1"), "mailto:garbageinfo@gmail.com", "", "", "garbageinfo@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:ayahya.testitg@gmail.com", "", "", "ayahya.testitg@gmail.com")

# Restore original (non-hyperlink-theme) formatting that Hyperlinks.Add overwrote
$ws.Range("A1").Font.Name = "Cambria"
$ws.Range("A1").Font.Underline = 2
$ws.Range("A1").Font.Color = 16711680
$ws.Range("A1").Font.Size = 11

$ws.Range("A2").Font.Name = "Arial"
$ws.Range("A2").Font.Underline = -4142
$ws.Range("A2").Font.Color = 0
$ws.Range("A2").Font.Size = 10

# Row 1's height grows slightly to fit the longer e-mail text
$ws.Rows("1").RowHeight = 14.15

$ws.Range("A1").Select()
